$d = $word.ActiveDocument

# 1) Month name: APRIL -> ABRIL (table "3. Data de Ocorrencia")
$d.Content.Find.Execute("APRIL", $true, $false, $false, $false, $false, $true, 1, $false, "ABRIL", 2) | Out-Null

# 2) Hour placeholder: ${oc_hr} -> 21:24
$d.Content.Find.Execute("`${oc_hr}", $true, $false, $false, $false, $false, $true, 1, $false, "21:24", 2) | Out-Null

# 3) "Mortos" count: 30 -> 10 (Danos Humanos table)
$d.Content.Find.Execute("30", $true, $false, $false, $false, $false, $true, 1, $false, "10", 2) | Out-Null

# 4) "Feridos" count: 0 -> 20 (Danos Humanos table, row right below "Mortos").
#    "0" is extremely common across the document, so Find/Replace can't be scoped
#    safely here (this runtime's Find operates document-wide regardless of the
#    range it's invoked on). Address the exact merged cell directly instead.
$humanTable = $d.Tables.Item(4)
$humanTable.Cell(4, 3).Range.Text = "20"

# 5) Valor (R$) for "Unidades Habitacionais": 1500.83 -> 150083.0
$d.Content.Find.Execute("1500.83", $true, $false, $false, $false, $false, $true, 1, $false, "150083.0", 2) | Out-Null
